$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.5
$ws.Range("C2").Value = 0.3846153846153846
$ws.Range("D2").Value = 0.4347826086956522
$ws.Range("E2").Value = 13
$ws.Range("B3").Value = 0.4285714285714285
$ws.Range("C3").Value = 0.5454545454545454
$ws.Range("D3").Value = 0.4799999999999999
$ws.Range("E3").Value = 11
$ws.Range("B4").Value = 0.4583333333333333
$ws.Range("C4").Value = 0.4583333333333333
$ws.Range("D4").Value = 0.4583333333333333
$ws.Range("E4").Value = 0.4583333333333333
$ws.Range("B5").Value = 0.4642857142857143
$ws.Range("C5").Value = 0.465034965034965
$ws.Range("D5").Value = 0.457391304347826
$ws.Range("E5").Value = 24
$ws.Range("B6").Value = 0.4672619047619048
$ws.Range("C6").Value = 0.4583333333333333
$ws.Range("D6").Value = 0.4555072463768116
$ws.Range("E6").Value = 24
$ws.Range("B7").Value = 0.5714285714285714
$ws.Range("C7").Value = 0.6153846153846154
$ws.Range("D7").Value = 0.5925925925925927
$ws.Range("E7").Value = 13
$ws.Range("B8").Value = 0.5
$ws.Range("C8").Value = 0.4545454545454545
$ws.Range("D8").Value = 0.4761904761904762
$ws.Range("E8").Value = 11
$ws.Range("B9").Value = 0.5416666666666666
$ws.Range("C9").Value = 0.5416666666666666
$ws.Range("D9").Value = 0.5416666666666666
$ws.Range("E9").Value = 0.5416666666666666
$ws.Range("B10").Value = 0.5357142857142857
$ws.Range("C10").Value = 0.534965034965035
$ws.Range("D10").Value = 0.5343915343915344
$ws.Range("E10").Value = 24
$ws.Range("B11").Value = 0.5386904761904762
$ws.Range("C11").Value = 0.5416666666666666
$ws.Range("D11").Value = 0.539241622574956
$ws.Range("E11").Value = 24
$ws.Range("B12").Value = 0.5384615384615384
$ws.Range("C12").Value = 0.5384615384615384
$ws.Range("D12").Value = 0.5384615384615384
$ws.Range("E12").Value = 13
$ws.Range("B13").Value = 0.4545454545454545
$ws.Range("C13").Value = 0.4545454545454545
$ws.Range("D13").Value = 0.4545454545454545
$ws.Range("E13").Value = 11
$ws.Range("B14").Value = 0.5
$ws.Range("C14").Value = 0.5
$ws.Range("D14").Value = 0.5
$ws.Range("E14").Value = 0.5
$ws.Range("B15").Value = 0.4965034965034965
$ws.Range("C15").Value = 0.4965034965034965
$ws.Range("D15").Value = 0.4965034965034965
$ws.Range("E15").Value = 24
$ws.Range("B16").Value = 0.5
$ws.Range("C16").Value = 0.5
$ws.Range("D16").Value = 0.5
$ws.Range("E16").Value = 24
$ws.Range("B17").Value = 0.6153846153846154
$ws.Range("C17").Value = 0.6153846153846154
$ws.Range("D17").Value = 0.6153846153846154
$ws.Range("E17").Value = 13
$ws.Range("B18").Value = 0.5454545454545454
$ws.Range("C18").Value = 0.5454545454545454
$ws.Range("D18").Value = 0.5454545454545454
$ws.Range("E18").Value = 11
$ws.Range("B19").Value = 0.5833333333333334
$ws.Range("C19").Value = 0.5833333333333334
$ws.Range("D19").Value = 0.5833333333333334
$ws.Range("E19").Value = 0.5833333333333334
$ws.Range("B20").Value = 0.5804195804195804
$ws.Range("C20").Value = 0.5804195804195804
$ws.Range("D20").Value = 0.5804195804195804
$ws.Range("E20").Value = 24
$ws.Range("B21").Value = 0.5833333333333334
$ws.Range("C21").Value = 0.5833333333333334
$ws.Range("D21").Value = 0.5833333333333334
$ws.Range("E21").Value = 24
$ws.Range("B22").Value = 0.5
$ws.Range("C22").Value = 0.3076923076923077
$ws.Range("D22").Value = 0.380952380952381
$ws.Range("E22").Value = 13
$ws.Range("B23").Value = 0.4375
$ws.Range("C23").Value = 0.6363636363636364
$ws.Range("D23").Value = 0.5185185185185185
$ws.Range("E23").Value = 11
$ws.Range("B24").Value = 0.4583333333333333
$ws.Range("C24").Value = 0.4583333333333333
$ws.Range("D24").Value = 0.4583333333333333
$ws.Range("E24").Value = 0.4583333333333333
$ws.Range("B25").Value = 0.46875
$ws.Range("C25").Value = 0.472027972027972
$ws.Range("D25").Value = 0.4497354497354498
$ws.Range("E25").Value = 24
$ws.Range("B26").Value = 0.4713541666666667
$ws.Range("C26").Value = 0.4583333333333333
$ws.Range("D26").Value = 0.4440035273368606
$ws.Range("E26").Value = 24
